$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.212.62"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.740.44"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").Value = "3.734.58"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.487"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "4.360.48"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "3.737.55"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "70.160.14"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "506.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "13.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000135"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +20.90%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.138"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.339"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("E40").Value = "  -6.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.39"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "431.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "3.012.73"
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.61"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("E51").Value = "  +1.87%  "
